$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the recorded timings in column B, rows 16-25
$values = @(90241, 75516, 80047, 71363, 74006, 107610, 114407, 72873, 89487, 74383)
$row = 16
foreach ($v in $values) {
    $ws.Cells.Item($row, 2).Value = $v
    $row++
}

# Average formula in B26
$ws.Range("B26").Formula = "=AVERAGE(B16:B25)"

# Update the selection to reflect where the user left off
$ws.Range("F24").Select()
